$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before B (old B/"vWrite" col and old C/"vExpected Result"
# col shift right to C/D); the new column inherits column A's formatting.
$ws.Columns("B:B").Insert()

# Re-purpose the (now 4-column) test-data template: A/B become the new
# "vName"/"vLastName" fields, C keeps the old "vWrite" header slot but is
# renamed to "vAlertText", and D keeps the original "vExpected Result" /
# "Searched with success" pair untouched.
$ws.Range("A1").Value = "vName"
$ws.Range("A2").Value = "'Matheus"
$ws.Range("B1").Value = "vLastName"
$ws.Range("B2").Value = "'Dias"
$ws.Range("C1").Value = "vAlertText"
$ws.Range("C2").Value = "Era Teste?"

$null = $ws.Range("C6").Select()
